$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Marking right/wrong values corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12: Total right/wrong values corrected, plus displayed fraction text
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "54 / 112"
